$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 <- original row 30
$ws.Range("A29").Value = 130671357
$ws.Range("B29").Value = 79243
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 6425
$ws.Range("F29").Value = "Garnlav"
$ws.Range("G29").Value = "Alectoria sarmentosa"
$ws.Range("H29").Value = "(Ach.) Ach."
$ws.Range("M29").ClearContents()
$ws.Range("Q29").Value = 557304
$ws.Range("R29").Value = 6710306
$ws.Range("Z29").Value = "09:53"
$ws.Range("AB29").Value = "09:53"

# Row 30 <- original row 29
$ws.Range("A30").Value = 130671332
$ws.Range("B30").Value = 91808
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 1202
$ws.Range("F30").Value = "Ullticka"
$ws.Range("G30").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H30").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("M30").ClearContents()
$ws.Range("Q30").Value = 556970
$ws.Range("R30").Value = 6710400
$ws.Range("Z30").Value = "11:06"
$ws.Range("AB30").Value = "11:06"

# Row 31 <- original row 33
$ws.Range("A31").Value = 130671363
$ws.Range("B31").Value = 79243
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 6425
$ws.Range("F31").Value = "Garnlav"
$ws.Range("G31").Value = "Alectoria sarmentosa"
$ws.Range("H31").Value = "(Ach.) Ach."
$ws.Range("M31").ClearContents()
$ws.Range("Q31").Value = 557265
$ws.Range("R31").Value = 6710358
$ws.Range("Z31").Value = "10:04"
$ws.Range("AB31").Value = "10:04"

# Row 33 <- original row 31
$ws.Range("A33").Value = 130671359
$ws.Range("B33").Value = 79243
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 6425
$ws.Range("F33").Value = "Garnlav"
$ws.Range("G33").Value = "Alectoria sarmentosa"
$ws.Range("H33").Value = "(Ach.) Ach."
$ws.Range("M33").ClearContents()
$ws.Range("Q33").Value = 557278
$ws.Range("R33").Value = 6710338
$ws.Range("Z33").Value = "10:00"
$ws.Range("AB33").Value = "10:00"

# Row 34 <- original row 37
$ws.Range("A34").Value = 130671327
$ws.Range("B34").Value = 91771
$ws.Range("D34").Value = "LC"
$ws.Range("E34").Value = 5447
$ws.Range("F34").Value = "Vedticka"
$ws.Range("G34").Value = "Fuscoporia viticola"
$ws.Range("H34").Value = "(Schwein.) Murrill"
$ws.Range("M34").ClearContents()
$ws.Range("Q34").Value = 556951
$ws.Range("R34").Value = 6710395
$ws.Range("Z34").Value = "11:01"
$ws.Range("AB34").Value = "11:01"

# Row 35 <- original row 34
$ws.Range("A35").Value = 130671362
$ws.Range("B35").Value = 79243
$ws.Range("D35").Value = "NT"
$ws.Range("E35").Value = 6425
$ws.Range("F35").Value = "Garnlav"
$ws.Range("G35").Value = "Alectoria sarmentosa"
$ws.Range("H35").Value = "(Ach.) Ach."
$ws.Range("M35").ClearContents()
$ws.Range("Q35").Value = 557271
$ws.Range("R35").Value = 6710350
$ws.Range("Z35").Value = "10:03"
$ws.Range("AB35").Value = "10:03"

# Row 37 <- original row 35
$ws.Range("A37").Value = 130671328
$ws.Range("B37").Value = 92179
$ws.Range("D37").Value = "VU"
$ws.Range("E37").Value = 2062
$ws.Range("F37").Value = "Ulltickeporing"
$ws.Range("G37").Value = "Skeletocutis brevispora"
$ws.Range("H37").Value = "Niemelä"
$ws.Range("M37").ClearContents()
$ws.Range("Q37").Value = 556973
$ws.Range("R37").Value = 6710401
$ws.Range("Z37").Value = "11:06"
$ws.Range("AB37").Value = "11:06"

# Row 39 <- original row 41
$ws.Range("A39").Value = 130671342
$ws.Range("B39").Value = 57881
$ws.Range("D39").Value = "NT"
$ws.Range("E39").Value = 100049
$ws.Range("F39").Value = "Spillkråka"
$ws.Range("G39").Value = "Dryocopus martius"
$ws.Range("H39").Value = "(Linnaeus, 1758)"
$ws.Range("M39").Value = "färska spår"
$ws.Range("Q39").Value = 557118
$ws.Range("R39").Value = 6710289
$ws.Range("Z39").Value = "11:27"
$ws.Range("AB39").Value = "11:27"

# Row 40 <- original row 39
$ws.Range("A40").Value = 130671338
$ws.Range("B40").Value = 57881
$ws.Range("D40").Value = "NT"
$ws.Range("E40").Value = 100049
$ws.Range("F40").Value = "Spillkråka"
$ws.Range("G40").Value = "Dryocopus martius"
$ws.Range("H40").Value = "(Linnaeus, 1758)"
$ws.Range("M40").Value = "gammalt bo"
$ws.Range("Q40").Value = 557146
$ws.Range("R40").Value = 6710445
$ws.Range("Z40").Value = "10:38"
$ws.Range("AB40").Value = "10:38"

# Row 41 <- original row 40
$ws.Range("A41").Value = 130671339
$ws.Range("B41").Value = 57881
$ws.Range("D41").Value = "NT"
$ws.Range("E41").Value = 100049
$ws.Range("F41").Value = "Spillkråka"
$ws.Range("G41").Value = "Dryocopus martius"
$ws.Range("H41").Value = "(Linnaeus, 1758)"
$ws.Range("M41").Value = "gammalt bo"
$ws.Range("Q41").Value = 557178
$ws.Range("R41").Value = 6710135
$ws.Range("Z41").Value = "11:44"
$ws.Range("AB41").Value = "11:44"

# Row 43 <- original row 44
$ws.Range("A43").Value = 130671361
$ws.Range("B43").Value = 79243
$ws.Range("D43").Value = "NT"
$ws.Range("E43").Value = 6425
$ws.Range("F43").Value = "Garnlav"
$ws.Range("G43").Value = "Alectoria sarmentosa"
$ws.Range("H43").Value = "(Ach.) Ach."
$ws.Range("M43").ClearContents()
$ws.Range("Q43").Value = 557273
$ws.Range("R43").Value = 6710349
$ws.Range("Z43").Value = "10:03"
$ws.Range("AB43").Value = "10:03"

# Row 44 <- original row 43
$ws.Range("A44").Value = 130671336
$ws.Range("B44").Value = 57881
$ws.Range("D44").Value = "NT"
$ws.Range("E44").Value = 100049
$ws.Range("F44").Value = "Spillkråka"
$ws.Range("G44").Value = "Dryocopus martius"
$ws.Range("H44").Value = "(Linnaeus, 1758)"
$ws.Range("M44").Value = "äldre spår"
$ws.Range("Q44").Value = 556993
$ws.Range("R44").Value = 6710383
$ws.Range("Z44").Value = "11:11"
$ws.Range("AB44").Value = "11:11"

# Row 45 <- original row 46
$ws.Range("A45").Value = 130671355
$ws.Range("B45").Value = 79243
$ws.Range("D45").Value = "NT"
$ws.Range("E45").Value = 6425
$ws.Range("F45").Value = "Garnlav"
$ws.Range("G45").Value = "Alectoria sarmentosa"
$ws.Range("H45").Value = "(Ach.) Ach."
$ws.Range("M45").ClearContents()
$ws.Range("Q45").Value = 557310
$ws.Range("R45").Value = 6710293
$ws.Range("Z45").Value = "09:52"
$ws.Range("AB45").Value = "09:52"

# Row 46 <- original row 45
$ws.Range("A46").Value = 130671335
$ws.Range("B46").Value = 91829
$ws.Range("D46").Value = "NT"
$ws.Range("E46").Value = 5442
$ws.Range("F46").Value = "Tallticka"
$ws.Range("G46").Value = "Porodaedalea pini"
$ws.Range("H46").Value = "(Brot.) Murrill"
$ws.Range("M46").ClearContents()
$ws.Range("Q46").Value = 557382
$ws.Range("R46").Value = 6710266
$ws.Range("Z46").Value = "12:03"
$ws.Range("AB46").Value = "12:03"

# Row 52 <- original row 53
$ws.Range("A52").Value = 130671330
$ws.Range("B52").Value = 91808
$ws.Range("D52").Value = "NT"
$ws.Range("E52").Value = 1202
$ws.Range("F52").Value = "Ullticka"
$ws.Range("G52").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H52").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("M52").ClearContents()
$ws.Range("Q52").Value = 557102
$ws.Range("R52").Value = 6710446
$ws.Range("Z52").Value = "10:34"
$ws.Range("AB52").Value = "10:34"

# Row 53 <- original row 52
$ws.Range("A53").Value = 130671325
$ws.Range("B53").Value = 5197
$ws.Range("D53").Value = "LC"
$ws.Range("E53").Value = 105930
$ws.Range("F53").Value = "Vågbandad barkbock"
$ws.Range("G53").Value = "Semanotus undatus"
$ws.Range("H53").Value = "(Linnaeus, 1758)"
$ws.Range("M53").Value = "färska gnagspår"
$ws.Range("Q53").Value = 557029
$ws.Range("R53").Value = 6710414
$ws.Range("Z53").Value = "10:51"
$ws.Range("AB53").Value = "10:51"

# Row 55 <- original row 57
$ws.Range("A55").Value = 130671358
$ws.Range("B55").Value = 79243
$ws.Range("D55").Value = "NT"
$ws.Range("E55").Value = 6425
$ws.Range("F55").Value = "Garnlav"
$ws.Range("G55").Value = "Alectoria sarmentosa"
$ws.Range("H55").Value = "(Ach.) Ach."
$ws.Range("M55").ClearContents()
$ws.Range("Q55").Value = 557300
$ws.Range("R55").Value = 6710292
$ws.Range("Z55").Value = "09:54"
$ws.Range("AB55").Value = "09:54"

# Row 56 <- original row 55
$ws.Range("A56").Value = 130671373
$ws.Range("B56").Value = 5177
$ws.Range("D56").Value = "LC"
$ws.Range("E56").Value = 100526
$ws.Range("F56").Value = "Bronshjon"
$ws.Range("G56").Value = "Callidium coriaceum"
$ws.Range("H56").Value = "Paykull, 1800"
$ws.Range("M56").Value = "färska gnagspår"
$ws.Range("Q56").Value = 557193
$ws.Range("R56").Value = 6710075
$ws.Range("Z56").Value = "11:51"
$ws.Range("AB56").Value = "11:51"

# Row 57 <- original row 56
$ws.Range("A57").Value = 130671329
$ws.Range("B57").Value = 91808
$ws.Range("D57").Value = "NT"
$ws.Range("E57").Value = 1202
$ws.Range("F57").Value = "Ullticka"
$ws.Range("G57").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H57").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("M57").ClearContents()
$ws.Range("Q57").Value = 557322
$ws.Range("R57").Value = 6710273
$ws.Range("Z57").Value = "09:49"
$ws.Range("AB57").Value = "09:49"
